{"js": "// Replace the date and each \"AxB=\" multiplication prompt in the table\n// with the new values from the target revision. Every original text\n// value in this document is unique, so a simple exact-text search and\n// replace (one search per pair) is unambiguous and safe regardless of\n// the order the pairs are applied in.\nconst replacements = [\n  [\"2024-07-27 Saturday\", \"2024-07-28 Sunday\"],\n  [\"90\u00d790=\", \"35\u00d784=\"],\n  [\"37\u00d762=\", \"49\u00d713=\"],\n  [\"21\u00d796=\", \"96\u00d712=\"],\n  [\"77\u00d767=\", \"51\u00d749=\"],\n  [\"23\u00d787=\", \"66\u00d781=\"],\n  [\"31\u00d744=\", \"43\u00d797=\"],\n  [\"16\u00d768=\", \"55\u00d736=\"],\n  [\"80\u00d796=\", \"18\u00d786=\"],\n  [\"79\u00d742=\", \"42\u00d724=\"],\n  [\"11\u00d743=\", \"83\u00d739=\"],\n  [\"97\u00d780=\", \"76\u00d780=\"],\n  [\"33\u00d740=\", \"95\u00d783=\"],\n  [\"35\u00d779=\", \"93\u00d740=\"],\n  [\"26\u00d736=\", \"83\u00d734=\"],\n  [\"60\u00d769=\", \"40\u00d719=\"],\n  [\"39\u00d746=\", \"87\u00d792=\"],\n  [\"20\u00d711=\", \"81\u00d797=\"],\n  [\"39\u00d751=\", \"43\u00d782=\"],\n  [\"75\u00d732=\", \"50\u00d714=\"],\n  [\"80\u00d732=\", \"31\u00d746=\"],\n  [\"45\u00d797=\", \"17\u00d744=\"],\n  [\"44\u00d772=\", \"36\u00d729=\"],\n  [\"44\u00d727=\", \"13\u00d798=\"],\n  [\"44\u00d776=\", \"83\u00d741=\"],\n  [\"66\u00d775=\", \"88\u00d732=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and each \"AxB=\" multiplication prompt in the table\n# with the new values from the target revision. Every original text\n# value in this document is unique, so a simple exact-text Find/Replace\n# (one pass per pair, ReplaceAll) is unambiguous and safe regardless of\n# the order the pairs are applied in.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-07-27 Saturday\", \"2024-07-28 Sunday\"),\n  @(\"90\u00d790=\", \"35\u00d784=\"),\n  @(\"37\u00d762=\", \"49\u00d713=\"),\n  @(\"21\u00d796=\", \"96\u00d712=\"),\n  @(\"77\u00d767=\", \"51\u00d749=\"),\n  @(\"23\u00d787=\", \"66\u00d781=\"),\n  @(\"31\u00d744=\", \"43\u00d797=\"),\n  @(\"16\u00d768=\", \"55\u00d736=\"),\n  @(\"80\u00d796=\", \"18\u00d786=\"),\n  @(\"79\u00d742=\", \"42\u00d724=\"),\n  @(\"11\u00d743=\", \"83\u00d739=\"),\n  @(\"97\u00d780=\", \"76\u00d780=\"),\n  @(\"33\u00d740=\", \"95\u00d783=\"),\n  @(\"35\u00d779=\", \"93\u00d740=\"),\n  @(\"26\u00d736=\", \"83\u00d734=\"),\n  @(\"60\u00d769=\", \"40\u00d719=\"),\n  @(\"39\u00d746=\", \"87\u00d792=\"),\n  @(\"20\u00d711=\", \"81\u00d797=\"),\n  @(\"39\u00d751=\", \"43\u00d782=\"),\n  @(\"75\u00d732=\", \"50\u00d714=\"),\n  @(\"80\u00d732=\", \"31\u00d746=\"),\n  @(\"45\u00d797=\", \"17\u00d744=\"),\n  @(\"44\u00d772=\", \"36\u00d729=\"),\n  @(\"44\u00d727=\", \"13\u00d798=\"),\n  @(\"44\u00d776=\", \"83\u00d741=\"),\n  @(\"66\u00d775=\", \"88\u00d732=\")\n)\n\nforeach ($p in $pairs) {\n  $find = $d.Content.Find\n  $find.Text = $p[0]\n  $find.Replacement.Text = $p[1]\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
